# Add season-record columns (Wins / Losses / Ties) to the DET_2002 sheet.
#
# The existing data goes through column AC. We append three new columns:
#   AD -> Wins
#   AE -> Losses
#   AF -> Ties
# The header row (row 1) reuses the same bold/centered/bordered header
# style already applied to the other header cells (e.g. AC1), and every
# data row (2-59) gets the team's season record: 55 wins, 106 losses,
# 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (from AC1) onto the new header
# cells so AD1:AF1 end up sharing the same style as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$lastRow = 59
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 55
    $ws.Cells.Item($r, 31).Value = 106
    $ws.Cells.Item($r, 32).Value = 0
}
